$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.06109198740482865
$ws.Cells.Item(2, 8).Value = -4.992384176756656
$ws.Cells.Item(2, 9).Value = -18.30683187187097
$ws.Cells.Item(3, 7).Value = 0.067676696442485
$ws.Cells.Item(3, 8).Value = 20.406085498212
$ws.Cells.Item(4, 7).Value = -0.02057860590377215
$ws.Cells.Item(4, 8).Value = -1.314508715363747
$ws.Cells.Item(5, 7).Value = -0.01870724044074673
$ws.Cells.Item(5, 8).Value = -63.64883931014639
$ws.Cells.Item(6, 7).Value = -0.02505877085537056
$ws.Cells.Item(6, 8).Value = -123.8677960677459
$ws.Cells.Item(7, 7).Value = -0.0233242581144179
$ws.Cells.Item(7, 8).Value = -311.8175364863462
$ws.Cells.Item(8, 7).Value = -0.003748757949902887
$ws.Cells.Item(8, 8).Value = 34.10703908896081
$ws.Cells.Item(9, 7).Value = -0.00008026014476906685
$ws.Cells.Item(9, 8).Value = 98.53890391227173
$ws.Cells.Item(10, 7).Value = -0.06291384425344568
$ws.Cells.Item(10, 8).Value = 0.1712272214593091
$ws.Cells.Item(11, 7).Value = -0.05993904288554885
$ws.Cells.Item(11, 8).Value = 6.459318334738437
$ws.Cells.Item(12, 7).Value = -0.3963079336387677
$ws.Cells.Item(12, 8).Value = -0.4319609984324115
$ws.Cells.Item(13, 7).Value = -0.422691146751019
$ws.Cells.Item(13, 8).Value = -7.796595167930128
$ws.Cells.Item(14, 7).Value = -0.01688172189131988
$ws.Cells.Item(14, 8).Value = -106.0187820311739
$ws.Cells.Item(15, 7).Value = -0.01304068755358358
$ws.Cells.Item(15, 8).Value = 71.22447330667818
$ws.Cells.Item(16, 7).Value = 0.1447753379616255
$ws.Cells.Item(16, 8).Value = 5.868891557688986
$ws.Cells.Item(17, 7).Value = 0.1439244079370158
$ws.Cells.Item(17, 8).Value = 3.193338374579437
$ws.Cells.Item(18, 7).Value = 0.1172087943696962
$ws.Cells.Item(18, 8).Value = -0.3904312487090952
$ws.Cells.Item(19, 7).Value = 0.1193034598386365
$ws.Cells.Item(19, 8).Value = -7.278678022929432
$ws.Cells.Item(20, 7).Value = 0.09160659898137542
$ws.Cells.Item(20, 8).Value = 3.23657877328293
$ws.Cells.Item(21, 7).Value = 0.09018216569223217
$ws.Cells.Item(21, 8).Value = 3.549321431609647
$ws.Cells.Item(22, 7).Value = -0.09182477211001129
$ws.Cells.Item(22, 8).Value = 1.774006280764877
$ws.Cells.Item(23, 7).Value = -0.1024413551517454
$ws.Cells.Item(23, 8).Value = -0.9821448562508609
$ws.Cells.Item(24, 7).Value = 0.1608108847432235
$ws.Cells.Item(24, 8).Value = -0.1739244507092207
$ws.Cells.Item(25, 7).Value = 0.1701397034994159
$ws.Cells.Item(25, 8).Value = -0.2673152829320175
$ws.Cells.Item(26, 7).Value = 0.08701313342394233
$ws.Cells.Item(26, 8).Value = -4.016006431138591
$ws.Cells.Item(27, 7).Value = 0.08254698651311029
$ws.Cells.Item(27, 8).Value = -3.972417865238016
$ws.Cells.Item(28, 7).Value = -0.1384982976314527
$ws.Cells.Item(28, 8).Value = -0.6232625331731505
$ws.Cells.Item(29, 7).Value = -0.1439768029557718
$ws.Cells.Item(29, 8).Value = -2.987871598373936
$ws.Cells.Item(30, 7).Value = 0.05156404773301053
$ws.Cells.Item(30, 8).Value = -0.8674891368363835
$ws.Cells.Item(31, 7).Value = 0.04714520345484524
$ws.Cells.Item(31, 8).Value = 7.58978948416903
$ws.Cells.Item(32, 7).Value = 0.1093347653979893
$ws.Cells.Item(32, 8).Value = 0.580322830128719
$ws.Cells.Item(33, 7).Value = 0.11895348956238
$ws.Cells.Item(33, 8).Value = -4.14335380322381
$ws.Cells.Item(34, 7).Value = -0.01201637087133283
$ws.Cells.Item(34, 8).Value = 23.05284490048182
$ws.Cells.Item(35, 7).Value = -0.01871457632773259
$ws.Cells.Item(35, 8).Value = -11.81366940981246
$ws.Cells.Item(36, 7).Value = 0.0317348054078177
$ws.Cells.Item(36, 8).Value = -13.68650283497494
$ws.Cells.Item(37, 7).Value = 0.03689980585449666
$ws.Cells.Item(37, 8).Value = 3.402368861327424
$ws.Cells.Item(38, 7).Value = 0.1001395117177247
$ws.Cells.Item(38, 8).Value = -0.1622273584623784
$ws.Cells.Item(39, 7).Value = 0.1024774038693842
$ws.Cells.Item(39, 8).Value = 5.212924148233674
$ws.Cells.Item(40, 7).Value = 0.04258022130560988
$ws.Cells.Item(40, 8).Value = 26.3972391795508
$ws.Cells.Item(41, 7).Value = 0.04228916665046199
$ws.Cells.Item(41, 8).Value = 31.25648704401194
$ws.Cells.Item(42, 7).Value = 0.1182671872457674
$ws.Cells.Item(42, 8).Value = -2.182071747362758
$ws.Cells.Item(43, 7).Value = 0.1217143384146356
$ws.Cells.Item(43, 8).Value = -4.751017371784244
$ws.Cells.Item(44, 7).Value = 0.03961865258102002
$ws.Cells.Item(44, 8).Value = -0.1209747943401801
$ws.Cells.Item(45, 7).Value = 0.04089023683135639
$ws.Cells.Item(45, 8).Value = 31.21678238360615
$ws.Cells.Item(46, 7).Value = 0.05951572564746652
$ws.Cells.Item(46, 8).Value = 5.11865093010923
$ws.Cells.Item(47, 7).Value = 0.06078580517578842
$ws.Cells.Item(47, 8).Value = 3.605951563294449
$ws.Cells.Item(48, 7).Value = 0.04908010696352166
$ws.Cells.Item(48, 8).Value = -0.3597754589338136
$ws.Cells.Item(49, 7).Value = 0.04660526343423394
$ws.Cells.Item(49, 8).Value = 2.276865141176867
$ws.Cells.Item(50, 7).Value = 0.02732719686470184
$ws.Cells.Item(50, 8).Value = 3.174121257390639
$ws.Cells.Item(51, 7).Value = 0.0258306650619133
$ws.Cells.Item(51, 8).Value = -7.799877420203607
$ws.Cells.Item(52, 7).Value = -0.08744666574696575
$ws.Cells.Item(52, 8).Value = -0.6017058525727982
$ws.Cells.Item(53, 7).Value = -0.08381474419893702
$ws.Cells.Item(53, 8).Value = -4.486054167769398
$ws.Cells.Item(54, 7).Value = 0.04185830603786439
$ws.Cells.Item(54, 8).Value = -16.32807195909839
$ws.Cells.Item(55, 7).Value = 0.05064042868611905
$ws.Cells.Item(55, 8).Value = -10.02864421669661
$ws.Cells.Item(56, 7).Value = 0.0523152772270862
$ws.Cells.Item(56, 8).Value = 5.833040065402835
$ws.Cells.Item(57, 7).Value = 0.04533536398827245
$ws.Cells.Item(57, 8).Value = 19.35391326725338
$ws.Cells.Item(58, 7).Value = 0.05502003170904462
$ws.Cells.Item(58, 8).Value = -4.502410060039372
$ws.Cells.Item(59, 7).Value = 0.05938952769526064
$ws.Cells.Item(59, 8).Value = 4.167831061078264
$ws.Cells.Item(60, 7).Value = 0.02798418471821184
$ws.Cells.Item(60, 8).Value = 1.908904621102251
$ws.Cells.Item(61, 7).Value = 0.02775485444626683
$ws.Cells.Item(61, 8).Value = 3.95923307198245
$ws.Cells.Item(62, 7).Value = 0.06145808582905458
$ws.Cells.Item(62, 8).Value = -1.583970053896864
$ws.Cells.Item(63, 7).Value = 0.06100690466595073
$ws.Cells.Item(63, 8).Value = -4.501635134083367
$ws.Cells.Item(64, 7).Value = 0.03210711929153812
$ws.Cells.Item(64, 8).Value = 15.74446937528394
$ws.Cells.Item(65, 7).Value = 0.03227725380273278
$ws.Cells.Item(65, 8).Value = -8.89145485408536
$ws.Cells.Item(66, 7).Value = 0.08893317158655797
$ws.Cells.Item(66, 8).Value = 14.48188687932543
$ws.Cells.Item(67, 7).Value = 0.07812510757941661
$ws.Cells.Item(67, 8).Value = -0.9455023383679962
$ws.Cells.Item(68, 7).Value = -0.02393700525425549
$ws.Cells.Item(68, 8).Value = -10.07699569601241
$ws.Cells.Item(69, 7).Value = -0.01631030257801391
$ws.Cells.Item(69, 8).Value = 14.79004986349637
$ws.Cells.Item(70, 7).Value = 0.07470937590687834
$ws.Cells.Item(70, 8).Value = 3.783399920933747
$ws.Cells.Item(71, 7).Value = 0.07575965578201758
$ws.Cells.Item(71, 8).Value = -4.610090932792088
$ws.Cells.Item(72, 7).Value = -0.1484221301178276
$ws.Cells.Item(72, 8).Value = 3.397098332618171
$ws.Cells.Item(73, 7).Value = -0.1455004648481245
$ws.Cells.Item(73, 8).Value = 4.948595479871373
$ws.Cells.Item(74, 7).Value = 0.1524565049374337
$ws.Cells.Item(74, 8).Value = 1.354422274815918
$ws.Cells.Item(75, 7).Value = 0.1530141787685479
$ws.Cells.Item(75, 8).Value = 1.70614814017952
$ws.Cells.Item(76, 7).Value = -0.01069541015616207
$ws.Cells.Item(76, 8).Value = -931.9156581066194
$ws.Cells.Item(77, 7).Value = -0.01019120307860143
$ws.Cells.Item(77, 8).Value = -361.5815689449836
$ws.Cells.Item(78, 7).Value = 0.09755909286077187
$ws.Cells.Item(78, 8).Value = 8.44929999084378
$ws.Cells.Item(79, 7).Value = 0.09991899633378908
$ws.Cells.Item(79, 8).Value = 3.114127047800431
$ws.Cells.Item(80, 7).Value = -0.218095592820294
$ws.Cells.Item(80, 8).Value = -0.7680077112054862
$ws.Cells.Item(81, 7).Value = -0.2144794231051746
$ws.Cells.Item(81, 8).Value = -0.640754945683107
$ws.Cells.Item(82, 7).Value = 0.1646579162420391
$ws.Cells.Item(82, 8).Value = -1.765775645437451
$ws.Cells.Item(83, 7).Value = 0.1770018760409023
$ws.Cells.Item(83, 8).Value = 0.5527181493017058
$ws.Cells.Item(84, 7).Value = 0.1110816330396745
$ws.Cells.Item(84, 8).Value = 4.665928730379957
$ws.Cells.Item(85, 7).Value = 0.1132343022493931
$ws.Cells.Item(85, 8).Value = 8.297240186070542
